$wb = $excel.ActiveWorkbook

# --- Reorder sheets: move APP_PROFILE before BD ---
$appProfile = $wb.Worksheets.Item("APP_PROFILE")
$bd = $wb.Worksheets.Item("BD")
$appProfile.Move($bd)

# --- Rename the trailing blank sheet into the new PHYSICAL_DOMAIN resource sheet ---
$phys = $wb.Worksheets.Item("Sheet1")
$phys.Name = "PHYSICAL_DOMAIN"

# --- Column widths matching the new resource table layout ---
$phys.Columns.Item(1).ColumnWidth = 20.296875
$phys.Columns.Item(2).ColumnWidth = 23.69921875
$phys.Columns.Item(3).ColumnWidth = 18.796875
$phys.Columns.Item(4).ColumnWidth = 20.69921875
$phys.Columns.Item(5).ColumnWidth = 24.59765625

# --- Populate data; write order controls shared-string table insertion order ---
$phys.Range("A2").Value = "physical_domain"
$phys.Range("B1").Value = "physical_domain_name"
$phys.Range("B2").Value = "markacidomain"
$phys.Range("A1").Value = "type"
$phys.Range("C1").Value = "vlan_pool_name"
$phys.Range("C2").Value = "mark_pool"

# --- Header row styling (bold, matches the other resource sheets) ---
$phys.Range("A1:E1").Font.Bold = $true

# --- Selection / active-sheet bookkeeping ---
$vlan = $wb.Worksheets.Item("VLAN_POOL")
$vlan.Activate()
$vlan.Range("A1:E1").Select()

$phys.Activate()
$phys.Range("C2").Select()
